# Update "dados bibi" metrics in metricas_retencao_anual.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: cohort 2021, period_index 5 -> num_customers 28 -> 29 (cohort_size 2654 unchanged)
$ws.Range("C22").Value = 29
$ws.Range("E22").Value = 29 / 2654

# Row 36: cohort 2024, period_index 1 -> num_customers 129 -> 130 (cohort_size 1930 unchanged)
$ws.Range("C36").Value = 130
$ws.Range("E36").Value = 130 / 1930

# Row 37: cohort 2025, period_index 0 -> num_customers and cohort_size 798 -> 808 (retention_rate stays 1)
$ws.Range("C37").Value = 808
$ws.Range("D37").Value = 808
$ws.Range("E37").Value = 1
